$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.313390374183655
$ws.Range("B1").Value = 1.573299288749695
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.041088819503784
$ws.Range("E1").Value = 0.8318157196044922
